# Read from txt, insights support added
# The expense log was re-imported; rows 20-38 now hold the data that used
# to live one row below (a leading "misc" entry on 29.11.24 was dropped
# during the re-import), the "Rent" category was renamed to "Rent & Bills",
# and the four trailing duplicate/garbage rows (39-42) were removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 20-38 with their new (shifted + recategorised) contents ---
# Force column A to plain text first so day/month-ambiguous dates like
# "1.12.24" aren't auto-converted into date serials by Excel's parser.
$ws.Range("A20:A38").NumberFormat = "@"

$ws.Range("A20").Value = "30.11.24"
$ws.Range("B20").Value = 25
$ws.Range("C20").Value = "pohe"
$ws.Range("D20").Value = "Food & Necessities"

$ws.Range("A21").Value = "30.11.24"
$ws.Range("B21").Value = 473
$ws.Range("C21").Value = "healthy things grocery"
$ws.Range("D21").Value = "Personal Care"

$ws.Range("A22").Value = "30.11.24"
$ws.Range("B22").Value = 80
$ws.Range("C22").Value = "paneer and veggie"
$ws.Range("D22").Value = "Food & Necessities"

$ws.Range("A23").Value = "1.12.24"
$ws.Range("B23").Value = 300
$ws.Range("C23").Value = "petrol"
$ws.Range("D23").Value = "Transportation"

$ws.Range("A24").Value = "1.12.24"
$ws.Range("B24").Value = 30
$ws.Range("C24").Value = "pinapple juice from piyush, gore 🔴🔴"
$ws.Range("D24").Value = "Miscellaneous"

$ws.Range("A25").Value = "1.12.24"
$ws.Range("B25").Value = 166
$ws.Range("C25").Value = "dinner"
$ws.Range("D25").Value = "Food & Necessities"

$ws.Range("A26").Value = "1.12.24"
$ws.Range("B26").Value = 790
$ws.Range("C26").Value = "harmosa eats"
$ws.Range("D26").Value = "Food & Necessities"

$ws.Range("A27").Value = "1.12.24"
$ws.Range("B27").Value = 450
$ws.Range("C27").Value = "maid"
$ws.Range("D27").Value = "Rent & Bills"

$ws.Range("A28").Value = "1.12.24"
$ws.Range("B28").Value = 6250
$ws.Range("C28").Value = "rent"
$ws.Range("D28").Value = "Rent & Bills"

$ws.Range("A29").Value = "2.12.24"
$ws.Range("B29").Value = 199
$ws.Range("C29").Value = "Netflix"
$ws.Range("D29").Value = "Entertainment"

$ws.Range("A30").Value = "2.12.24"
$ws.Range("B30").Value = 151
$ws.Range("C30").Value = "rice, aata"
$ws.Range("D30").Value = "Food & Necessities"

$ws.Range("A31").Value = "2.12.24"
$ws.Range("B31").Value = 150
$ws.Range("C31").Value = "dinner"
$ws.Range("D31").Value = "Food & Necessities"

$ws.Range("A32").Value = "3.12.24"
$ws.Range("B32").Value = 70
$ws.Range("C32").Value = "paratha"
$ws.Range("D32").Value = "Food & Necessities"

$ws.Range("A33").Value = "3.12.24"
$ws.Range("B33").Value = 93
$ws.Range("C33").Value = "dinner"
$ws.Range("D33").Value = "Food & Necessities"

$ws.Range("A34").Value = "4.12.24"
$ws.Range("B34").Value = 111
$ws.Range("C34").Value = "lunch"
$ws.Range("D34").Value = "Food & Necessities"

$ws.Range("A35").Value = "4.12.24"
$ws.Range("B35").Value = 141
$ws.Range("C35").Value = "dinner"
$ws.Range("D35").Value = "Food & Necessities"

$ws.Range("A36").Value = "6.12.24"
$ws.Range("B36").Value = 214
$ws.Range("C36").Value = "electricity bill"
$ws.Range("D36").Value = "Rent & Bills"

$ws.Range("A37").Value = "6.12.24"
$ws.Range("B37").Value = 38
$ws.Range("C37").Value = "cab for party"
$ws.Range("D37").Value = "Transportation"

$ws.Range("A38").Value = "7.12.24"
$ws.Range("B38").Value = 77
$ws.Range("C38").Value = "lunch"
$ws.Range("D38").Value = "Food & Necessities"

# --- Remove the now-obsolete trailing rows 39-42 ---
$ws.Range("A39:D42").ClearContents()
